$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-19"

# Update the October row label (row 11, column A) to reflect the new date
$ws.Range("A11").Value = "October (through 10-19)"

# Update October row (row 11) values
$ws.Range("B11").Value = 18
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 47
$ws.Range("F11").Value = 28
$ws.Range("G11").Value = 90
$ws.Range("H11").Value = 124

# Update Total row (row 12) values
$ws.Range("B12").Value = 244
$ws.Range("C12").Value = 459
$ws.Range("D12").Value = 659
$ws.Range("E12").Value = 595
$ws.Range("F12").Value = 450
$ws.Range("G12").Value = 991
$ws.Range("H12").Value = 1372
